$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting existing data rows (2-20) down to (3-21)
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above (header); clear it
# so it matches the plain (unstyled) data rows, then restore the date-column
# style (s="2") on D2 by copying just the format from D3.
$ws.Range("A2:R2").ClearFormats()
$ws.Range("D3").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row 2 with the latest weekly price observation
$ws.Cells.Item(2, 1).Value = 3
$ws.Cells.Item(2, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(2, 3).Value = "Coquimbo"
$ws.Cells.Item(2, 4).Value = 44537
$ws.Cells.Item(2, 5).Value = 5
$ws.Cells.Item(2, 6).Value = 100112044
$ws.Cells.Item(2, 7).Value = "Perejil"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 88
$ws.Cells.Item(2, 11).Value = 2000
$ws.Cells.Item(2, 12).Value = 2200
$ws.Cells.Item(2, 13).Value = 2091
$ws.Cells.Item(2, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(2, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(2, 16).Value = 697
$ws.Cells.Item(2, 17).Value = 3
$ws.Cells.Item(2, 18).Value = "Hortaliza"
